# The commit reshuffles the data rows (rows 2-9) of the "Artfynd" sheet:
# each destination row's full contents come from a different source row
# (a permutation of rows 2..9), with no other structural changes.
#
# Destination row -> source row mapping (derived from the OOXML diff):
#   2 <- 3
#   3 <- 6
#   4 <- 8
#   5 <- 9
#   6 <- 2
#   7 <- 4
#   8 <- 5
#   9 <- 7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every source row's full value array (A:AY) BEFORE any writes,
# since several rows both give and receive data (the mapping has cycles).
$row2 = $ws.Range("A2:AY2").Value()
$row3 = $ws.Range("A3:AY3").Value()
$row4 = $ws.Range("A4:AY4").Value()
$row5 = $ws.Range("A5:AY5").Value()
$row6 = $ws.Range("A6:AY6").Value()
$row7 = $ws.Range("A7:AY7").Value()
$row8 = $ws.Range("A8:AY8").Value()
$row9 = $ws.Range("A9:AY9").Value()

# Columns Y/Z/AA/AB hold dates/times stored as plain text (e.g. "2022-04-13"),
# and column I ("Antal") holds small counts stored as plain text (e.g. "1").
# Force text format on those columns first so the re-assignment below doesn't
# let Excel auto-convert the number/date-looking strings into real numbers
# or date serials.
$ws.Range("Y2:AB9").NumberFormat = "@"
$ws.Range("I2:I9").NumberFormat = "@"

# Now write each destination row from the snapshot taken above.
$ws.Range("A2:AY2").Value() = $row3
$ws.Range("A3:AY3").Value() = $row6
$ws.Range("A4:AY4").Value() = $row8
$ws.Range("A5:AY5").Value() = $row9
$ws.Range("A6:AY6").Value() = $row2
$ws.Range("A7:AY7").Value() = $row4
$ws.Range("A8:AY8").Value() = $row5
$ws.Range("A9:AY9").Value() = $row7
